$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57 (pushes old rows 57-82 down to 58-83),
# mirroring the data of the (old) row 57 but with a new "Fecha" (date) value.
$ws.Rows.Item(57).Insert()

$ws.Range("A57").Value2 = 9
$ws.Range("B57").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C57").Value2 = "Metropolitana"
$ws.Range("D57").Value2 = 45146
$ws.Range("E57").Value2 = 13
$ws.Range("F57").Value2 = 100112010
$ws.Range("G57").Value2 = "Achicoria"
$ws.Range("H57").Value2 = "Sin especificar"
$ws.Range("I57").Value2 = "Primera"
$ws.Range("J57").Value2 = 70
$ws.Range("K57").Value2 = 7000
$ws.Range("L57").Value2 = 7000
$ws.Range("M57").Value2 = 7000
$ws.Range("N57").Value2 = "$/caja 16 unidades"
$ws.Range("O57").Value2 = "Provincia de Quillota"
$ws.Range("P57").Value2 = 438
$ws.Range("Q57").Value2 = 16
$ws.Range("R57").Value2 = "Hortaliza"
